$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, pushing existing rows 42:54 down to 43:55
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new weekly record
$ws.Range("A42").Value = 5
$ws.Range("B42").Value = "Macroferia Regional de Talca"
$ws.Range("C42").Value = "Maule"
$ws.Range("D42").Value = 44504
$ws.Range("E42").Value = 7
$ws.Range("F42").Value = 100112022
$ws.Range("G42").Value = "Arveja Verde"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 600
$ws.Range("K42").Value = 12000
$ws.Range("L42").Value = 12000
$ws.Range("M42").Value = 12000
$ws.Range("N42").Value = "$/saco 25 kilos"
$ws.Range("O42").Value = "Región del Maule"
$ws.Range("P42").Value = 480
$ws.Range("Q42").Value = 25
$ws.Range("R42").Value = "Hortaliza"
